# honda.xlsx: "drop platform overrides"
#
# The "test" sheet had a duplicate/extra row for "AndroidChromeTest2"
# (row 9) that is removed entirely, and the explicit "local" platform
# override in column B is cleared for the AndroidChromeTest and
# AndroidNativeTest rows (they now fall back to whatever default applies
# instead of forcing "local").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test")

# Remove the whole "AndroidChromeTest2" row; everything below shifts up
# a row (so the former row 10/11/12 become row 9/10/11).
$ws.Rows("9:9").Delete()

# Drop the platform overrides ("local") for AndroidChromeTest (row 8)
# and AndroidNativeTest (now row 9 after the delete above).
$ws.Range("B8").ClearContents()
$ws.Range("B9").ClearContents()

# Reset the lingering selection back to the default top-left cell.
$ws.Range("A1").Select()
